# Add a new "2022-Q1" worksheet with fund-level holding data, positioned
# right before the "总计" (totals) sheet, and update the "总计" sheet with
# a new summary row for 2022-Q1.
#
# NOTE: worksheet references returned by Worksheets.Item(...) are
# resolved positionally under the hood, so any variable captured before
# a sheet is inserted/moved/deleted can silently start pointing at the
# wrong sheet afterwards. To stay safe, sheets are re-fetched by name
# immediately before each use, right after any structural change.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet by cloning the most similarly
#    shaped existing quarter sheet ("2021-Q4": header + bold/bordered
#    index column + plain data cells), which carries over sheetPr,
#    styles and dimension automatically. It has one extra data row, so
#    we drop it afterwards to go from 4 data rows down to 3.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q4").Copy($wb.Worksheets.Item("总计"))

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"
$newSheet = $wb.Worksheets.Item("2022-Q1")
$newSheet.Rows.Item(5).Delete()
$newSheet = $wb.Worksheets.Item("2022-Q1")

# A cell known to carry no special style (on an unrelated, untouched
# sheet), used below to strip the formatting that gets implicitly
# applied when NumberFormat is set on the fund-data cells.
$cleanCell = $wb.Worksheets.Item("2021-Q3").Range("B2")

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Force the fund-code / fund-name / numeric-looking columns to be
# stored as text (matching the source data, e.g. "005536" keeps its
# leading zero and "11.34" isn't turned into a float).
$newSheet.Range("B2:G4").NumberFormat = "@"

# Row 2 - 广发资管平衡精选一年持有混合A
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "870009"
$newSheet.Cells.Item(2,3).Value = "广发资管平衡精选一年持有混合A"
$newSheet.Cells.Item(2,4).Value = "11.34"
$newSheet.Cells.Item(2,5).Value = "94.29"
$newSheet.Cells.Item(2,6).Value = "9.12"
$newSheet.Cells.Item(2,7).Value = "1.0342"
$newSheet.Cells.Item(2,8).Value = 4

# Row 3 - 广发资管平衡精选一年持有混合C
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "872019"
$newSheet.Cells.Item(3,3).Value = "广发资管平衡精选一年持有混合C"
$newSheet.Cells.Item(3,4).Value = "1.54"
$newSheet.Cells.Item(3,5).Value = "94.29"
$newSheet.Cells.Item(3,6).Value = "9.12"
$newSheet.Cells.Item(3,7).Value = "0.1404"
$newSheet.Cells.Item(3,8).Value = 4

# Row 4 - 渤海汇金量化成长混合
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "005536"
$newSheet.Cells.Item(4,3).Value = "渤海汇金量化成长混合"
$newSheet.Cells.Item(4,4).Value = "0.61"
$newSheet.Cells.Item(4,5).Value = "88.57"
$newSheet.Cells.Item(4,6).Value = "0.77"
$newSheet.Cells.Item(4,7).Value = "0.0047"
$newSheet.Cells.Item(4,8).Value = 6

# Setting NumberFormat="@" leaves a style footprint (numFmtId 49) on the
# text cells; strip it back out by pasting in the "clean" (style-less)
# formatting from the known unstyled cell, which keeps the cell
# contents (still text) but drops the extraneous style index, matching
# the unstyled B:G data cells used throughout the workbook.
$cleanCell.Copy()
$newSheet.Range("B2:G4").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new first data row for 2022-Q1
#    and shift the existing rows' index column down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet = $wb.Worksheets.Item("总计")

# Restore formatting on the inserted row: column A should match the
# other bold/bordered index cells, while B:D should be plain/unstyled
# (Insert() otherwise carries the header's bold style onto the whole row).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial($xlPasteFormats)
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial($xlPasteFormats)

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 3
$totalSheet.Cells.Item(2,4).Value = 1.18

# Renumber the index column for the rows that were pushed down
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(5,1).Value = 3
